# Delete the "employee_id" column (column A) from the employees sheet.
# This shifts name/position/salary one column to the left (A/B/C),
# matching the source data's new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A:A").Delete()

# Mirror the cursor position recorded after the edit.
$ws.Range("C10").Select()
